# Apply updated NATMI TPM-based ligand-receptor statistics (Kitl-Kit sheet)
# G:J = ligand expression/specificity stats (vary by sending cluster)
# M:T = receptor expression/specificity + edge weight stats (vary by target cluster)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    "G2", "H2", "I2", "J2", "M2", "N2", "O2", "P2", "Q2", "R2",
    "S2", "T2", "G3", "H3", "I3", "J3", "M3", "O3", "P3", "Q3",
    "R3", "S3", "T3", "G4", "H4", "I4", "J4", "M4", "N4", "O4",
    "P4", "Q4", "R4", "S4", "T4", "G5", "H5", "I5", "J5", "M5",
    "N5", "O5", "P5", "Q5", "R5", "S5", "T5", "I6", "J6", "M6",
    "N6", "O6", "P6", "Q6", "R6", "S6", "T6", "I7", "J7", "M7",
    "O7", "P7", "Q7", "R7", "S7", "T7", "I8", "J8", "M8", "N8",
    "O8", "P8", "Q8", "R8", "S8", "T8", "I9", "J9", "M9", "N9",
    "O9", "P9", "Q9", "R9", "S9", "T9", "G10", "H10", "I10", "J10",
    "M10", "N10", "O10", "P10", "Q10", "R10", "S10", "T10", "G11", "H11",
    "I11", "J11", "M11", "O11", "P11", "Q11", "R11", "S11", "T11", "G12",
    "H12", "I12", "J12", "M12", "N12", "O12", "P12", "Q12", "R12", "S12",
    "T12", "G13", "H13", "I13", "J13", "M13", "N13", "O13", "P13", "Q13",
    "R13", "S13", "T13", "G14", "H14", "I14", "J14", "M14", "N14", "O14",
    "P14", "Q14", "R14", "S14", "T14", "G15", "H15", "I15", "J15", "M15",
    "O15", "P15", "Q15", "R15", "S15", "T15", "G16", "H16", "I16", "J16",
    "M16", "N16", "O16", "P16", "Q16", "R16", "S16", "T16", "G17", "H17",
    "I17", "J17", "M17", "N17", "O17", "P17", "Q17", "R17", "S17", "T17"
)

$values = @(
    147.0057066666667, 441.01712, 0.8587975787179232, 0.8587975787179231, 19.38942866666667, 58.16828600000001, 0.9489681539286383, 0.9489681539286383, 2850.356663006258, 25653.20996705632,
    0.814971552874332, 0.8149715528743319, 147.0057066666667, 441.01712, 0.8587975787179232, 0.8587975787179231, 0.01112833333333333, 0.000544649051871798, 0.000544649051871798, 1.635928505688889,
    14.7233565512, 0.0004677432869985127, 0.0004677432869985126, 147.0057066666667, 441.01712, 0.8587975787179232, 0.8587975787179231, 0.1557643333333334, 0.4672930000000001, 0.007623504250301878,
    0.007623504250301877, 22.89824589512889, 206.08421305616, 0.006547046991505049, 0.006547046991505047, 147.0057066666667, 441.01712, 0.8587975787179232, 0.8587975787179231, 0.8757959999999999,
    2.627388, 0.04286369276918794, 0.04286369276918795, 128.74700987584, 1158.72308888256, 0.03681123556508756, 0.03681123556508756, 0.03365062030176851, 0.03365062030176851, 19.38942866666667,
    58.16828600000001, 0.9489681539286383, 0.9489681539286383, 111.6867026274462, 1005.180323647016, 0.03193336702632282, 0.03193336702632282, 0.03365062030176851, 0.03365062030176851, 0.01112833333333333,
    0.000544649051871798, 0.000544649051871798, 0.0641012624511111, 0.5769113620599998, 0.00001832777844225609, 0.00001832777844225609, 0.03365062030176851, 0.03365062030176851, 0.1557643333333334, 0.4672930000000001,
    0.007623504250301878, 0.007623504250301877, 0.8972314283231112, 8.075082854908, 0.0002565356468958268, 0.0002565356468958268, 0.03365062030176851, 0.03365062030176851, 0.8757959999999999, 2.627388,
    0.04286369276918794, 0.04286369276918795, 5.044747274191999, 45.40272546772799, 0.001442389850107604, 0.001442389850107604, 14.74162533333333, 44.224876, 0.08611959650886204, 0.08611959650886204,
    19.38942866666667, 58.16828600000001, 0.9489681539286383, 0.9489681539286383, 285.831692831393, 2572.485235482537, 0.08172475451609401, 0.08172475451609401, 14.74162533333333, 44.224876,
    0.08611959650886204, 0.08611959650886204, 0.01112833333333333, 0.000544649051871798, 0.000544649051871798, 0.1640497205844444, 1.47644748526, 0.00004690495658613352, 0.00004690495658613352, 14.74162533333333,
    44.224876, 0.08611959650886204, 0.08611959650886204, 0.1557643333333334, 0.4672930000000001, 0.007623504250301878, 0.007623504250301877, 2.296219442296445, 20.665974980668, 0.0006565331100195926,
    0.0006565331100195924, 14.74162533333333, 44.224876, 0.08611959650886204, 0.08611959650886204, 0.8757959999999999, 2.627388, 0.04286369276918794, 0.04286369276918795, 12.910656500432,
    116.195908503888, 0.003691403926162293, 0.003691403926162294, 3.668683333333334, 11.00605, 0.02143220447144637, 0.02143220447144637, 19.38942866666667, 58.16828600000001, 0.9489681539286383,
    0.9489681539286383, 71.13367379225558, 640.2030641303002, 0.02033847951188957, 0.02033847951188957, 3.668683333333334, 11.00605, 0.02143220447144637, 0.02143220447144637, 0.01112833333333333,
    0.000544649051871798, 0.000544649051871798, 0.04082633102777778, 0.36743697925, 0.00001167302984489577, 0.00001167302984489577, 3.668683333333334, 11.00605, 0.02143220447144637, 0.02143220447144637,
    0.1557643333333334, 0.4672930000000001, 0.007623504250301878, 0.007623504250301877, 0.571450013627778, 5.143050122650002, 0.0001633885018814103, 0.0001633885018814103, 3.668683333333334, 11.00605,
    0.02143220447144637, 0.02143220447144637, 0.8757959999999999, 2.627388, 0.04286369276918794, 0.04286369276918795, 3.2130181886, 28.9171636974, 0.0009186634278304932, 0.0009186634278304934
)

for ($i = 0; $i -lt $cells.Length; $i++) {
    $ws.Range($cells[$i]).Value = $values[$i]
}
